$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.210.98"
$ws.Range("E2").Value = "  +1.59%  "

$ws.Range("D3").Value = "2.246.40"
$ws.Range("E3").Value = "  +1.19%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "'306.49"
$ws.Range("E5").Value = "  -1.87%  "

$ws.Range("D6").Value = "'96.06"
$ws.Range("E6").Value = "  -1.38%  "

$ws.Range("E7").Value = "  +1.13%  "

$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  -0.74%  "

$ws.Range("D10").Value = "'35.24"
$ws.Range("E10").Value = "  -1.24%  "

$ws.Range("D11").Value = "'0.0816"
$ws.Range("E11").Value = "  -0.55%  "

$ws.Range("E12").Value = "  -1.47%  "

$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").Value = "2.589.43"
$ws.Range("E14").Value = "  +1.00%  "

$ws.Range("D15").Value = "2.324.75"
$ws.Range("E15").Value = "  +4.32%  "

$ws.Range("D16").Value = "'0.835"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").Value = "'13.63"
$ws.Range("E17").Value = "  -2.96%  "

$ws.Range("D18").Value = "43.934.01"
$ws.Range("E18").Value = "  +1.10%  "

$ws.Range("D19").Value = "0.0₃0971"
$ws.Range("E19").Value = "  +0.69%  "

$ws.Range("D20").Value = "'6.40"
$ws.Range("E20").Value = "  +2.08%  "

$ws.Range("D21").Value = "'12.16"
$ws.Range("E21").Value = "  -5.96%  "

$ws.Range("D22").Value = "'65.49"
$ws.Range("E22").Value = "  +0.51%  "

$ws.Range("D23").Value = "'237.15"
$ws.Range("E23").Value = "  +1.16%  "

$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("E25").Value = "  -1.29%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").Value = "'10.00"
$ws.Range("E27").Value = "  +0.23%  "

$ws.Range("D28").Value = "'2.21"
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").Value = "'37.72"
$ws.Range("E29").Value = "  +4.39%  "

$ws.Range("D30").Value = "'6.00"
$ws.Range("E30").Value = "  +1.28%  "

$ws.Range("D31").Value = "'20.18"
$ws.Range("E31").Value = "  +1.85%  "

$ws.Range("D32").Value = "'152.94"
$ws.Range("E32").Value = "  -4.85%  "

$ws.Range("D33").Value = "'0.0802"
$ws.Range("E33").Value = "  -2.66%  "

$ws.Range("E34").Value = "  +4.35%  "

$ws.Range("D35").Value = "'2.61"
$ws.Range("E35").Value = "  -2.88%  "

$ws.Range("D36").Value = "'0.120"
$ws.Range("E36").Value = "  +3.39%  "

$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("E38").Value = "  -5.49%  "

$ws.Range("E39").Value = "  -1.39%  "

$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'14.68"
$ws.Range("E40").Value = "  -4.27%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'3.87"
$ws.Range("E41").Value = "  -2.04%  "

$ws.Range("E42").Value = "  -2.51%  "

$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  -0.43%  "

$ws.Range("D44").Value = "1.744.18"
$ws.Range("E44").Value = "  +2.46%  "

$ws.Range("D45").Value = "'83.15"
$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("E46").Value = "  -0.79%  "

$ws.Range("D47").Value = "'100.29"
$ws.Range("E47").Value = "  -0.93%  "

$ws.Range("D48").Value = "'4.96"
$ws.Range("E48").Value = "  -2.54%  "

$ws.Range("D49").Value = "'8.16"
$ws.Range("E49").Value = "  +2.08%  "

$ws.Range("D50").Value = "'54.99"
$ws.Range("E50").Value = "  -2.24%  "

$ws.Range("E51").Value = "  -3.68%  "
